# Adds cache for moves generator: new shared string, and a new block of
# benchmark rows (26-29) on the first worksheet, mirroring the existing
# rows 18-20 / 22-24 "run" blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Copy cell formatting (number format / borders / fill / font) from
#    existing analogous cells so that the workbook's style table is
#    reused exactly as Excel itself would dedupe it, instead of minting
#    brand-new style records.
# ---------------------------------------------------------------------

# Row 26 mirrors row 22 (start-of-run row).
$ws.Range("A22:F22").Copy() | Out-Null
$ws.Range("A26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("H22").Copy() | Out-Null
$ws.Range("H26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("I22:L22").Copy() | Out-Null
$ws.Range("I26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("N22").Copy() | Out-Null
$ws.Range("N26").PasteSpecial($xlPasteFormats) | Out-Null

# G26/M26 use the "Bad" (red) style already present on M22 (bordered),
# since both deltas below turn out negative.
$ws.Range("M22").Copy() | Out-Null
$ws.Range("G26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M22").Copy() | Out-Null
$ws.Range("M26").PasteSpecial($xlPasteFormats) | Out-Null

# Row 27 mirrors row 23.
$ws.Range("C23:L23").Copy() | Out-Null
$ws.Range("C27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("N23").Copy() | Out-Null
$ws.Range("N27").PasteSpecial($xlPasteFormats) | Out-Null
# M27's delta is positive -> reuse the "Good" (un-bordered) style, as
# already used on G23.
$ws.Range("G23").Copy() | Out-Null
$ws.Range("M27").PasteSpecial($xlPasteFormats) | Out-Null

# Row 28 mirrors row 24.
$ws.Range("I24:L24").Copy() | Out-Null
$ws.Range("I28").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("N24").Copy() | Out-Null
$ws.Range("N28").PasteSpecial($xlPasteFormats) | Out-Null
# M28's delta is positive -> same "Good" style as M27/G23.
$ws.Range("G23").Copy() | Out-Null
$ws.Range("M28").PasteSpecial($xlPasteFormats) | Out-Null

# Row 29 is a brand new depth level (7) appended at the end of the run;
# I/J/K/L follow the plain pattern of row 24, while N29 reuses the
# bordered "Bad" style already present on M23/M24 (no corresponding M29
# cell exists, since there's no prior-run comparison row).
$ws.Range("J24").Copy() | Out-Null
$ws.Range("J29").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L24").Copy() | Out-Null
$ws.Range("L29").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M23").Copy() | Out-Null
$ws.Range("N29").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in values / formulas for the new rows.
# ---------------------------------------------------------------------

# Row 26
$ws.Range("A26").Value = 45563
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 206603
$ws.Range("E26").Value = 619
$ws.Range("F26").Formula = "=D26/E26*1000"
$ws.Range("G26").Formula = "=(E22-E26)/E22"
$ws.Range("H26").Formula = "=(F26-80000000)/80000000"
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 197281
$ws.Range("K26").Value = 25
$ws.Range("L26").Formula = "=J26/K26*1000"
$ws.Range("M26").Formula = "=(K22-K26)/K22"
$ws.Range("N26").Formula = "=(L26-80000000)/80000000"
$ws.Range("P26").Value = "moves cache"

# Row 27
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 5072212
$ws.Range("E27").Value = 11632
$ws.Range("F27").Formula = "=D27/E27*1000"
$ws.Range("G27").Formula = "=(E23-E27)/E23"
$ws.Range("H27").Formula = "=(F27-80000000)/80000000"
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 4880523
$ws.Range("K27").Value = 429
$ws.Range("L27").Formula = "=J27/K27*1000"
$ws.Range("M27").Formula = "=(K23-K27)/K23"
$ws.Range("N27").Formula = "=(L27-80000000)/80000000"

# Row 28
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 119060324
$ws.Range("K28").Value = 7892
$ws.Range("L28").Formula = "=J28/K28*1000"
$ws.Range("M28").Formula = "=(K24-K28)/K24"
$ws.Range("N28").Formula = "=(L28-80000000)/80000000"

# Row 29
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 3195901860
$ws.Range("K29").Value = 151609
$ws.Range("L29").Formula = "=J29/K29*1000"
$ws.Range("N29").Formula = "=(L29-80000000)/80000000"

# ---------------------------------------------------------------------
# 3. Misc sheet-level bookkeeping to mirror the author's final state.
# ---------------------------------------------------------------------

$ws.Range("L29").Select() | Out-Null

Write-Host "edit applied"
